$d = $word.ActiveDocument

# Change 1: reframe the opening sentence about service contribution and add
# the new "Beyond the standard..." sentence introducing John's service.
$d.Content.Find.Execute(
    "community truly deserves recognition. He advocates",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "research community deserves recognition. Beyond the standard of reviewing papers, serving on panels, and committee work, John advocates",
    2) | Out-Null

# Change 2: add the "[See CV section on service]" citation and change
# "is a role model" to "sets the bar".
$d.Content.Find.Execute(
    "for all to work in." + [char]160 + " John is a role model in how",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "for all to work in [See CV section on service]." + [char]160 + " John sets the bar in how",
    2) | Out-Null

# Change 3: rewrite the presentation title and add an additional year (2120)
# to the AGU meeting citation.
$d.Content.Find.Execute(
    [char]8220 + "How Using Team Science Ensured Safe Space-Time Travel" + [char]8221 + " Plenary at the AGU meeting 1920, and 2020]",
    $true, $false, $false, $false, $false, $true, 1, $false,
    [char]8220 + "It Takes a Team to Innovate New Physics and Safe Space-Time Travel" + [char]8221 + ", Plenary at the AGU meeting 1920, 2020, 2120]",
    2) | Out-Null

# Change 4: lower-case the "Collaborator" at the end of the paragraph.
$d.Content.Find.Execute(
    "good Collaborator.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "good collaborator.",
    2) | Out-Null

Write-Host "Done"
